# NBB_synthetic_GDP.xlsx -- update NACE 64 labels in all economic calibration data
#
# The header row (row 1) holds short NACE-Rev.2 industry codes (e.g. "01", "05-09",
# "64", ...). This change prefixes every one of those numeric/range codes with its
# NACE section letter (A, B, C, ... T), e.g. "01" -> "A01", "64" -> "K64".
# The "date" header (A1) and the "BE" header (BM1, Belgium country code) are left
# untouched since they are not NACE codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the NACE code labels in row 1 (B1:BL1) -----------------------------
$ws.Range("B1").Value  = "A01"
$ws.Range("C1").Value  = "A02"
$ws.Range("D1").Value  = "A03"
$ws.Range("E1").Value  = "B05-09"
$ws.Range("F1").Value  = "C10-12"
$ws.Range("G1").Value  = "C13-15"
$ws.Range("H1").Value  = "C16"
$ws.Range("I1").Value  = "C17"
$ws.Range("J1").Value  = "C18"
$ws.Range("K1").Value  = "C19"
$ws.Range("L1").Value  = "C20"
$ws.Range("M1").Value  = "C21"
$ws.Range("N1").Value  = "C22"
$ws.Range("O1").Value  = "C23"
$ws.Range("P1").Value  = "C24"
$ws.Range("Q1").Value  = "C25"
$ws.Range("R1").Value  = "C26"
$ws.Range("S1").Value  = "C27"
$ws.Range("T1").Value  = "C28"
$ws.Range("U1").Value  = "C29"
$ws.Range("V1").Value  = "C30"
$ws.Range("W1").Value  = "C31-32"
$ws.Range("X1").Value  = "C33"
$ws.Range("Y1").Value  = "D35"
$ws.Range("Z1").Value  = "E36"
$ws.Range("AA1").Value = "E37-39"
$ws.Range("AB1").Value = "F41-43"
$ws.Range("AC1").Value = "G45"
$ws.Range("AD1").Value = "G46"
$ws.Range("AE1").Value = "G47"
$ws.Range("AF1").Value = "H49"
$ws.Range("AG1").Value = "H50"
$ws.Range("AH1").Value = "H51"
$ws.Range("AI1").Value = "H52"
$ws.Range("AJ1").Value = "H53"
$ws.Range("AK1").Value = "I55-56"
$ws.Range("AL1").Value = "J58"
$ws.Range("AM1").Value = "J59-60"
$ws.Range("AN1").Value = "J61"
$ws.Range("AO1").Value = "J62-63"
$ws.Range("AP1").Value = "K64"
$ws.Range("AQ1").Value = "K65"
$ws.Range("AR1").Value = "K66"
$ws.Range("AS1").Value = "L68"
$ws.Range("AT1").Value = "M69-70"
$ws.Range("AU1").Value = "M71"
$ws.Range("AV1").Value = "M72"
$ws.Range("AW1").Value = "M73"
$ws.Range("AX1").Value = "M74-75"
$ws.Range("AY1").Value = "N77"
$ws.Range("AZ1").Value = "N78"
$ws.Range("BA1").Value = "N79"
$ws.Range("BB1").Value = "N80-82"
$ws.Range("BC1").Value = "O84"
$ws.Range("BD1").Value = "P85"
$ws.Range("BE1").Value = "Q86"
$ws.Range("BF1").Value = "Q87-88"
$ws.Range("BG1").Value = "R90-92"
$ws.Range("BH1").Value = "R93"
$ws.Range("BI1").Value = "S94"
$ws.Range("BJ1").Value = "S95"
$ws.Range("BK1").Value = "S96"
$ws.Range("BL1").Value = "T97-98"
# A1 ("date") and BM1 ("BE") keep their original text.

# --- Re-apply the header formatting that the original sheet used for these cells
# (numeric-style integer format, Calibri 11, centered) -------------------------
$headerRange = $ws.Range("B1:AE1")
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 11
$headerRange.Font.Color = 0
$headerRange.NumberFormat = "0"
$headerRange.HorizontalAlignment = -4108

# The taller font used above makes row 1 match the row height already used by
# every other row in the sheet.
$ws.Rows.Item(1).RowHeight = 13.8

# --- Restore the view/selection state ------------------------------------------
$ws.Range("I33").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
